$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Update timestamps for rows 2-12 in column A
$newTimestamp = "2025-10-16 06:35:01"
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Swap title (B) and URL (F) between row 6 and row 7
# (use Value2 for reads - Value getter is unreliable in this runtime)
$b6 = $ws.Cells.Item(6, 2).Value2
$b7 = $ws.Cells.Item(7, 2).Value2
$f6 = $ws.Cells.Item(6, 6).Value2
$f7 = $ws.Cells.Item(7, 6).Value2

$ws.Cells.Item(6, 2).Value = $b7
$ws.Cells.Item(7, 2).Value = $b6
$ws.Cells.Item(6, 6).Value = $f7
$ws.Cells.Item(7, 6).Value = $f6
